$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clone the last existing row (97) into new row 98 via copy/paste-special
# so that the text-like values (e.g. the "07.04.23" date-looking string)
# are copied verbatim as shared strings instead of being re-interpreted
# as dates by value-assignment type inference.
$ws.Range("A97:D97").Copy()
$ws.Range("A98:D98").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# This new row represents a "Create Country" test run, but deliberately
# marked as FAILED on purpose (per commit message) instead of PASSED.
$ws.Range("B98").Value = "FAILED"
